$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new trade row (row 4) mirroring the structure of row 3.
$ws.Range("A4").Value = 10193.879999999999
$ws.Range("B4").Value = 10119
$ws.Range("C4").Value = 20.3
$ws.Range("D4").Value = 20.149999999999999
$ws.Range("E4").Value = $true
$ws.Range("F4").Value = -0.74
$ws.Range("G4").Value = 42608.640474537038
$ws.Range("H4").Value = $true

# Match the date formatting used in column G (e.g. G3) for the new cell.
$ws.Range("G3").Copy()
$ws.Range("G4").PasteSpecial(-4122)
$ws.Range("G4").Value = 42608.640474537038

# Column A's "best fit" width grows now that it holds a wider value
# (longest entry goes from 5 to 8 significant characters).
$ws.Columns.Item(1).ColumnWidth = 8.14
